$wb = $excel.ActiveWorkbook

# New remn_amt values for row 100 and row 101 (col B), plus a brand-new
# row 102 (date 2025-10-31 / serial 45961, remn_amt 0) for each of the
# four worksheets in the workbook.
$updates = @(
    @{ Sheet = 1; B100 = 1224423; B101 = 1269642 },
    @{ Sheet = 2; B100 = 1100077; B101 = 1074719 },
    @{ Sheet = 3; B100 = 1793818; B101 = 1775506 },
    @{ Sheet = 4; B100 = 941867;  B101 = 1013785 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    # Update existing cells B100 / B101
    $ws.Range("B100").Value = $u.B100
    $ws.Range("B101").Value = $u.B101

    # Append new row 102
    $ws.Range("A102").Value = 45961
    $ws.Range("A102").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B102").Value = 0
}
